# Round 2 dense-answer-relevance results: add the D ("Answer_relevance")
# scores, wrap/size the Query & Response columns, and zoom the sheet in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths + wrap text for columns B (Query) and C (Response) ---
# Setting WrapText on the whole column picks up the existing bold/bordered
# header style in row 1 (B1/C1) and merges wrapText into it, while the
# plain body cells (B2:C71) get a separate wrap-only style - matching the
# two new cellXfs entries needed.
$ws.Columns.Item(2).ColumnWidth = 26.5
$ws.Columns.Item(3).ColumnWidth = 100
$ws.Columns.Item(2).WrapText = $true
$ws.Columns.Item(3).WrapText = $true

# --- Row 1 header gets a touch more height ---
$ws.Rows.Item(1).RowHeight = 16

# --- Per-row heights (rows 2-71), sized to fit the wrapped text ---
$rowHeights = @(320, 288, 304, 256, 350, 256, 240, 335, 224, 256, 272, 350, 256, 80, 304, 304, 256, 256, 335, 304, 176, 304, 288, 272, 272, 350, 240, 64, 304, 224, 320, 288, 365, 208, 272, 304, 240, 320, 272, 365, 288, 64, 304, 350, 256, 350, 320, 224, 395, 304, 272, 304, 256, 335, 304, 80, 304, 256, 320, 304, 365, 240, 256, 304, 240, 320, 256, 350, 240, 64)

# --- New column D values: whether the answer was judged relevant (1/0) ---
$dValues = @(1, 1, 1, 1, 0, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 0, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 0, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 0, 1, 1, 1, 1, 1, 1, 1, 1, 0, 1, 1, 1, 1, 0, 1, 1, 1, 1, 1, 1, 1, 1, 0)

for ($i = 0; $i -lt $rowHeights.Length; $i++) {
    $r = $i + 2
    $ws.Rows.Item($r).RowHeight = $rowHeights[$i]
    $ws.Cells.Item($r, 4).Value = $dValues[$i]
}

# --- View: zoom in to 140% and move the selection to D1 ---
$excel.ActiveWindow.Zoom = 140
[void]$ws.Range("D1").Select()
